$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D retains text formatting so numeric-looking strings
# (e.g. "1.00", "498.79") are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "54.189.45"
$ws.Range("E2").Value = "  -0.38%  "

# Row 3
$ws.Range("D3").Value = "2.272.85"
$ws.Range("E3").Value = "  +1.15%  "

# Row 4
$ws.Range("E4").Value = "  -0.56%  "

# Row 5
$ws.Range("D5").Value = "498.79"
$ws.Range("E5").Value = "  +0.76%  "

# Row 6
$ws.Range("D6").Value = "128.92"
$ws.Range("E6").Value = "  +1.18%  "

# Row 7
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.47%  "

# Row 8
$ws.Range("D8").Value = "0.526"
$ws.Range("E8").Value = "  -0.80%  "

# Row 9
$ws.Range("D9").Value = "0.0952"
$ws.Range("E9").Value = "  +0.13%  "

# Row 10
$ws.Range("D10").Value = "0.152"
$ws.Range("E10").Value = "  +0.47%  "

# Row 11
$ws.Range("D11").Value = "0.335"
$ws.Range("E11").Value = "  +3.33%  "

# Row 12
$ws.Range("D12").Value = "4.70"
$ws.Range("E12").Value = "  +1.42%  "

# Row 13
$ws.Range("D13").Value = "2.673.92"
$ws.Range("E13").Value = "  -0.81%  "

# Row 14
$ws.Range("D14").Value = "22.74"
$ws.Range("E14").Value = "  +4.82%  "

# Row 15
$ws.Range("D15").Value = "54.149.63"
$ws.Range("E15").Value = "  -0.43%  "

# Row 16
$ws.Range("E16").Value = "  -0.46%  "

# Row 17
$ws.Range("D17").Value = "2.282.26"
$ws.Range("E17").Value = "  -0.35%  "

# Row 18
$ws.Range("D18").Value = "10.20"
$ws.Range("E18").Value = "  +1.59%  "

# Row 19
$ws.Range("D19").Value = "4.14"
$ws.Range("E19").Value = "  +1.90%  "

# Row 20
$ws.Range("D20").Value = "302.13"
$ws.Range("E20").Value = "  -0.80%  "

# Row 21
$ws.Range("E21").Value = "  -2.49%  "

# Row 22
$ws.Range("E22").Value = "  -0.03%  "

# Row 23
$ws.Range("D23").Value = "61.52"
$ws.Range("E23").Value = "  -3.34%  "

# Row 24
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  -0.50%  "

# Row 25
$ws.Range("D25").Value = "0.149"
$ws.Range("E25").Value = "  -1.63%  "

# Row 26
$ws.Range("D26").Value = "7.27"
$ws.Range("E26").Value = "  +2.38%  "

# Row 27
$ws.Range("D27").Value = "170.43"
$ws.Range("E27").Value = "  +0.24%  "

# Row 28
$ws.Range("E28").Value = "  +0.19%  "

# Row 29
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").Value = "5.92"
$ws.Range("E29").Value = "  +0.56%  "

# Row 30
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0683"
$ws.Range("E30").Value = "  -0.85%  "

# Row 31
$ws.Range("E31").Value = "  +0.40%  "

# Row 32
$ws.Range("E32").Value = "  +0.30%  "

# Row 33
$ws.Range("D33").Value = "17.69"
$ws.Range("E33").Value = "  +0.23%  "

# Row 34
$ws.Range("D34").Value = "0.956"
$ws.Range("E34").Value = "  +9.87%  "

# Row 35
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  +0.64%  "

# Row 36
$ws.Range("E36").Value = "  -1.50%  "

# Row 37
$ws.Range("D37").Value = "3.69"
$ws.Range("E37").Value = "  +1.14%  "

# Row 38
$ws.Range("D38").Value = "0.372"
$ws.Range("E38").Value = "  -0.95%  "

# Row 39
$ws.Range("E39").Value = "  -0.27%  "

# Row 40
$ws.Range("D40").Value = "3.36"
$ws.Range("E40").Value = "  +0.48%  "

# Row 41
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "4.81"
$ws.Range("E41").Value = "  -2.50%  "

# Row 42
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "125.03"
$ws.Range("E42").Value = "  -3.93%  "

# Row 43
$ws.Range("D43").Value = "0.0493"
$ws.Range("E43").Value = "  +2.48%  "

# Row 44
$ws.Range("D44").Value = "0.0892"
$ws.Range("E44").Value = "  -0.31%  "

# Row 45
$ws.Range("D45").Value = "0.545"
$ws.Range("E45").Value = "  -0.61%  "

# Row 46
$ws.Range("D46").Value = "238.99"
$ws.Range("E46").Value = "  -1.15%  "

# Row 47
$ws.Range("D47").Value = "0.371"
$ws.Range("E47").Value = "  -0.98%  "

# Row 48
$ws.Range("D48").Value = "0.0204"
$ws.Range("E48").Value = "  +0.38%  "

# Row 49
$ws.Range("E49").Value = "  +0.56%  "

# Row 50
$ws.Range("D50").Value = "16.20"
$ws.Range("E50").Value = "  -0.89%  "

# Row 51
$ws.Range("D51").Value = "4.64"
$ws.Range("E51").Value = "  -0.27%  "
